$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as TEXT (not auto-converted to a number),
# matching the source workbook where these price/volume cells are
# stored as strings even when their content looks numeric (e.g. "238.86",
# or multi-dot values like "29.361.94" that are not valid Excel numbers).
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "29.361.94"
Set-TextValue "E2" "  +0.01%  "

# Row 3
Set-TextValue "D3" "1.841.78"
Set-TextValue "E3" "  -0.18%  "

# Row 4
Set-TextValue "D4" "0.9991"
Set-TextValue "E4" "  +0.07%  "

# Row 5
Set-TextValue "D5" "238.86"

# Row 6
Set-TextValue "D6" "0.6261"
Set-TextValue "E6" "  -0.20%  "

# Row 7
Set-TextValue "D7" "1.000"
Set-TextValue "E7" "  -0.02%  "

# Row 8
Set-TextValue "D8" "0.07401"
Set-TextValue "E8" "  -0.81%  "

# Row 9
Set-TextValue "D9" "0.2891"
Set-TextValue "E9" "  -0.18%  "

# Row 10
Set-TextValue "D10" "24.81"
Set-TextValue "E10" "  +1.42%  "

# Row 11
Set-TextValue "D11" "0.07710"
Set-TextValue "E11" "  -0.40%  "

# Row 12
Set-TextValue "D12" "1.831.65"
Set-TextValue "E12" "  -0.75%  "

# Row 13
Set-TextValue "D13" "4.966"
Set-TextValue "E13" "  -0.28%  "

# Row 14
Set-TextValue "D14" "0.6728"
Set-TextValue "E14" "  -0.92%  "

# Row 15
Set-TextValue "D15" "0.00001024"
Set-TextValue "E15" "  -1.92%  "

# Row 16
Set-TextValue "D16" "81.73"
Set-TextValue "E16" "  -0.23%  "

# Row 17
Set-TextValue "D17" "6.244"
Set-TextValue "E17" "  +0.76%  "

# Row 18
Set-TextValue "D18" "29.321.19"
Set-TextValue "E18" "  -0.25%  "

# Row 19
Set-TextValue "D19" "234.56"
Set-TextValue "E19" "  +2.81%  "

# Row 20
Set-TextValue "D20" "12.31"
Set-TextValue "E20" "  +0.04%  "

# Row 21
Set-TextValue "E21" "  +0.02%  "

# Row 22
Set-TextValue "D22" "7.301"
Set-TextValue "E22" "  -2.69%  "

# Row 24
Set-TextValue "D24" "158.00"
Set-TextValue "E24" "  -0.59%  "

# Row 25
Set-TextValue "D25" "8.479"
Set-TextValue "E25" "  +0.29%  "

# Row 26
Set-TextValue "D26" "0.1348"
Set-TextValue "E26" "  -1.31%  "

# Row 27
Set-TextValue "D27" "17.31"
Set-TextValue "E27" "  -1.13%  "

# Row 28
Set-TextValue "D28" "0.07220"
Set-TextValue "E28" "  +11.33%  "

# Row 29
Set-TextValue "D29" "1.478"
Set-TextValue "E29" "  +4.38%  "

# Row 30
Set-TextValue "D30" "1.474"
Set-TextValue "E30" "  -0.79%  "

# Row 31
Set-TextValue "D31" "4.057"
Set-TextValue "E31" "  -0.69%  "

# Row 32
Set-TextValue "D32" "4.024"
Set-TextValue "E32" "  -1.52%  "

# Row 33
Set-TextValue "D33" "1.818"
Set-TextValue "E33" "  -0.65%  "

# Row 34
Set-TextValue "D34" "1.142"
Set-TextValue "E34" "  +0.17%  "

# Row 35
Set-TextValue "D35" "0.6955"
Set-TextValue "E35" "  +0.22%  "

# Row 36
Set-TextValue "E36" "  -0.70%  "

# Row 37
Set-TextValue "D37" "6.915"
Set-TextValue "E37" "  +2.04%  "

# Row 38
Set-TextValue "D38" "0.01834"
Set-TextValue "E38" "  +0.03%  "

# Row 39
Set-TextValue "D39" "2.811"
Set-TextValue "E39" "  -0.76%  "

# Row 40
Set-TextValue "D40" "1.232.09"
Set-TextValue "E40" "  -2.32%  "

# Row 41
Set-TextValue "D41" "0.9436"
Set-TextValue "E41" "  +2.60%  "

# Row 42
Set-TextValue "E42" "  +0.03%  "

# Row 43
Set-TextValue "D43" "1.998.69"
Set-TextValue "E43" "  -0.21%  "

# Row 44
Set-TextValue "D44" "100.81"
Set-TextValue "E44" "  -0.72%  "

# Row 45
Set-TextValue "D45" "65.24"
Set-TextValue "E45" "  -1.22%  "

# Row 46
Set-TextValue "D46" "0.00000000120"
Set-TextValue "E46" "  +3.76%  "

# Row 47
Set-TextValue "D47" "1.703"
Set-TextValue "E47" "  -1.94%  "

# Row 48
Set-TextValue "D48" "6.943"
Set-TextValue "E48" "  -1.66%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D49" "0.3894"
Set-TextValue "E49" "  -1.30%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "8.813"
Set-TextValue "E50" "  -1.98%  "

# Row 51
Set-TextValue "D51" "0.1127"
Set-TextValue "E51" "  -2.50%  "

